# Sprint7.1Observaciones.xlsx - update per commit:
# "se hace mas grande el codigo de barras por que no lo alcanza a leer la pistola"
#
# Real user-visible edits captured by the diff:
#   1. Column B is resized narrower (was bestFit/162 chars -> ~26.29 chars, custom width).
#   2. B22 text tweak: "cambios" -> "cambio" (singular).
#   3. Two new rows of notes appended at the bottom of the sheet (APP / alert note),
#      with the selection left on the new empty cell below them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Narrow column B (engine snaps ColumnWidth to its own character grid; 25.5
#    is the closest input that lands nearest the target stored width).
$ws.Columns.Item(2).ColumnWidth = 25.5

# 2) Correct the wording in the existing ticket note (singular "cambio").
$ws.Range("B22").Value = "en ticket agrear cambio y con cuanto pago"

# 3) Append the new observation block.
$ws.Range("A33").Value = "APP"
$ws.Range("B34").Value = "Alerta cuando llega un pedido y no se atendido un notificacion y ademas un campanita que tiene  un pedio por atender"

# Leave the selection where the user ended up after typing the last entry.
$ws.Range("A34").Select() | Out-Null
